{"js": "// Replace each two-digit multiplication equation with its updated value.\n// Each \"old\" string is a unique, exact cell value, so Body.search with\n// matchCase + exact matching lets us target the precise <w:t> run and\n// swap its text in place, preserving all run/paragraph formatting.\nconst replacements = [\n  [\"83\u00d740=3320\", \"38\u00d742=1596\"],\n  [\"40\u00d771=2840\", \"55\u00d792=5060\"],\n  [\"20\u00d777=1540\", \"83\u00d754=4482\"],\n  [\"44\u00d794=4136\", \"35\u00d783=2905\"],\n  [\"74\u00d738=2812\", \"24\u00d722=528\"],\n  [\"92\u00d739=3588\", \"26\u00d735=910\"],\n  [\"82\u00d755=4510\", \"67\u00d735=2345\"],\n  [\"46\u00d725=1150\", \"45\u00d767=3015\"],\n  [\"78\u00d758=4524\", \"16\u00d738=608\"],\n  [\"25\u00d779=1975\", \"13\u00d759=767\"],\n  [\"93\u00d712=1116\", \"56\u00d738=2128\"],\n  [\"91\u00d768=6188\", \"77\u00d713=1001\"],\n  [\"14\u00d762=868\", \"80\u00d750=4000\"],\n  [\"62\u00d798=6076\", \"87\u00d770=6090\"],\n  [\"71\u00d734=2414\", \"80\u00d715=1200\"],\n  [\"34\u00d715=510\", \"88\u00d747=4136\"],\n  [\"28\u00d740=1120\", \"90\u00d786=7740\"],\n  [\"57\u00d799=5643\", \"94\u00d791=8554\"],\n  [\"64\u00d746=2944\", \"34\u00d741=1394\"],\n  [\"30\u00d766=1980\", \"17\u00d720=340\"],\n  [\"98\u00d762=6076\", \"16\u00d795=1520\"],\n  [\"47\u00d723=1081\", \"71\u00d749=3479\"],\n  [\"12\u00d796=1152\", \"75\u00d736=2700\"],\n  [\"61\u00d798=5978\", \"36\u00d746=1656\"],\n  [\"27\u00d762=1674\", \"62\u00d755=3410\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication equation with its updated value.\n# Every \"old\" string below is a unique, exact cell value in the table, so\n# Word's Find/Replace (case-sensitive, whole match) can safely retarget\n# each one without touching neighbouring cells or formatting.\n$pairs = @(\n  @(\"83\u00d740=3320\", \"38\u00d742=1596\"),\n  @(\"40\u00d771=2840\", \"55\u00d792=5060\"),\n  @(\"20\u00d777=1540\", \"83\u00d754=4482\"),\n  @(\"44\u00d794=4136\", \"35\u00d783=2905\"),\n  @(\"74\u00d738=2812\", \"24\u00d722=528\"),\n  @(\"92\u00d739=3588\", \"26\u00d735=910\"),\n  @(\"82\u00d755=4510\", \"67\u00d735=2345\"),\n  @(\"46\u00d725=1150\", \"45\u00d767=3015\"),\n  @(\"78\u00d758=4524\", \"16\u00d738=608\"),\n  @(\"25\u00d779=1975\", \"13\u00d759=767\"),\n  @(\"93\u00d712=1116\", \"56\u00d738=2128\"),\n  @(\"91\u00d768=6188\", \"77\u00d713=1001\"),\n  @(\"14\u00d762=868\", \"80\u00d750=4000\"),\n  @(\"62\u00d798=6076\", \"87\u00d770=6090\"),\n  @(\"71\u00d734=2414\", \"80\u00d715=1200\"),\n  @(\"34\u00d715=510\", \"88\u00d747=4136\"),\n  @(\"28\u00d740=1120\", \"90\u00d786=7740\"),\n  @(\"57\u00d799=5643\", \"94\u00d791=8554\"),\n  @(\"64\u00d746=2944\", \"34\u00d741=1394\"),\n  @(\"30\u00d766=1980\", \"17\u00d720=340\"),\n  @(\"98\u00d762=6076\", \"16\u00d795=1520\"),\n  @(\"47\u00d723=1081\", \"71\u00d749=3479\"),\n  @(\"12\u00d796=1152\", \"75\u00d736=2700\"),\n  @(\"61\u00d798=5978\", \"36\u00d746=1656\"),\n  @(\"27\u00d762=1674\", \"62\u00d755=3410\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n\n  $found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Could not find text to replace: $oldText\"\n  }\n}\n"}
